# Rewrite row 1 of the single sheet: the two existing "gamelog URL" /
# "final score" columns are replaced by a wider row that leads with
# Sims, Cam / WR identity columns, keeps the original game-log fields
# (shifted right), and appends a trailing numeric 0 column.
#
# A leading apostrophe forces Excel to store "2018-09-09", "1", "22.246"
# and the blank K1 cell as literal text instead of auto-converting them
# to a date / numbers; ClearFormats() afterwards drops the transient
# "stored as text" cell format so the cells end up back on the sheet's
# default style, matching plain untouched text cells like G1:J1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Sims"
$ws.Range("B1").Value = "Cam"
$ws.Range("C1").Value = "WR"
$ws.Range("D1").Value = "'2018-09-09"
$ws.Range("E1").Value = "'1"
$ws.Range("F1").Value = "'22.246"
$ws.Range("G1").Value = "WAS"
$ws.Range("H1").Value = "@"
$ws.Range("I1").Value = "ARI"
$ws.Range("J1").Value = "W 24-6"
$ws.Range("K1").Value = "'"
$ws.Range("L1").Value = 0

$ws.Range("D1:F1").ClearFormats()
$ws.Range("K1").ClearFormats()
